# Simulator full-month coverage: populate Rate and Total columns that were
# previously left at 0 for the "Weekly Timesheet" and "Jason Schema" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Timesheet" ---
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")

# Row 2: Hall, 8 hrs
$ws1.Range("E2").Value = 140
$ws1.Range("F2").Value = 1120

# Row 3: Hall, 9 hrs
$ws1.Range("E3").Value = 140
$ws1.Range("F3").Value = 1260

# Row 4: McGill, 7 hrs
$ws1.Range("E4").Value = 140
$ws1.Range("F4").Value = 980

# Subtotal / Hourly subtotal / Grand total rows
$ws1.Range("F6").Value = 3360
$ws1.Range("F10").Value = 3360
$ws1.Range("F11").Value = 3360

# --- Sheet 2: "Jason Schema" ---
$ws2 = $wb.Worksheets.Item("Jason Schema")

# Row 2: Hall, 8 hrs
$ws2.Range("F2").Value = 140
$ws2.Range("G2").Value = 1120

# Row 3: Hall, 9 hrs
$ws2.Range("F3").Value = 140
$ws2.Range("G3").Value = 1260

# Row 4: McGill, 7 hrs
$ws2.Range("F4").Value = 140
$ws2.Range("G4").Value = 980
